# "Code for Justify Exceptions"
# Rebuild the TestRun data grid: drop the unused "Date" column, rename the
# ID/description headers, replace the sample data row with a new failing
# test case (EmpID 10649101) and keep the second row (10510273) in sync
# with the column shift. Also re-point the conditional formatting and the
# _FilterDatabase defined name from column G to column F now that the
# sheet is one column narrower.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Grab "clean" donor formats before we touch anything, stashing them
#    far out in column J (well outside the A:G working range) so the
#    later column delete doesn't clobber them before we're done reading.
# ---------------------------------------------------------------------
$ws.Range("D2").Copy()                      # non-quote-prefixed date style (FromDate/ToDate)
$ws.Range("J1").PasteSpecial(-4122)

$ws.Range("A3").Copy()                      # clean EmpNum/EmpID style
$ws.Range("J2").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Remove column C ("Date") entirely - D:G shift left to become C:F.
#    The donor cells in J1/J2 shift left to I1/I2.
# ---------------------------------------------------------------------
$ws.Columns("C").Delete()

# ---------------------------------------------------------------------
# 3) Header row
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "EmpID"
$ws.Range("B1").Value = "Paycodes"

# ---------------------------------------------------------------------
# 4) Row 2 becomes the new sample/test row (EmpID 10649101). Values that
#    look numeric/date-like are entered quote-prefixed so Excel keeps
#    them as literal text, then the quote-prefix flag is wiped by
#    pasting the clean donor formats back on top (value is untouched by
#    a formats-only paste).
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "'10649101"
$ws.Range("I2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("B2").Value = "SK-Early Out Excused.,SK-Long Break Excused.,SK-Early Out Excused."

$ws.Range("C2").Value = "'01/02/2025"
$ws.Range("D2").Value = "'28/02/2025"
$ws.Range("I1").Copy()
$ws.Range("C2:D2").PasteSpecial(-4122)

$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "Pass"

# ---------------------------------------------------------------------
# 5) Row 3 (EmpID 10510273) already has the right values after the
#    column shift; just normalise the FromDate/ToDate cell style to the
#    same clean (non quote-prefixed) style used above, and drop the
#    stray TestResult value - the new layout leaves it blank.
# ---------------------------------------------------------------------
$ws.Range("I1").Copy()
$ws.Range("C3:D3").PasteSpecial(-4122)

$ws.Range("F3").ClearContents()

# ---------------------------------------------------------------------
# 6) Tidy up the scratch donor cells.
# ---------------------------------------------------------------------
$ws.Range("I1:I2").Clear()

# ---------------------------------------------------------------------
# 7) Conditional formatting followed column G; move it onto column F.
# ---------------------------------------------------------------------
$fcs = $ws.Range("G1:G1048576").FormatConditions
$fcs.Item(1).ModifyAppliesToRange($ws.Range("F1:F1048576"))

# ---------------------------------------------------------------------
# 8) The _FilterDatabase defined name also referenced column G.
# ---------------------------------------------------------------------
$name = $wb.Names.Item(1)
$name.RefersTo = "='Full Time 40'!`$F`$1:`$F`$3"

# ---------------------------------------------------------------------
# 9) Restore the expected selection.
# ---------------------------------------------------------------------
$ws.Range("B2").Select()
